$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.328.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.937.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.93%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7217'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.93%  '

$ws.Range("E7").Value = '  -0.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3314'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.21'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07271'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8064'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08072'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.935.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.492'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.330.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008213'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.828'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.189.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.22%  '

$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("E23").Value = '  -1.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.991'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.739'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.362'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1314'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.567'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.90%  '

$ws.Range("E31").Value = '  -2.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.415'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.161'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05181'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.270'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7454'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.745'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01970'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.814'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.403'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4520'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.014'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8476'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("E45").Value = '  -0.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.81%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.446'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.62%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.680'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4182'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06038'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.65%  '
